$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 9-17 (data/tickers that are no longer part of the input set)
$ws.Range("A9:C17").EntireRow.Delete()

# New ticker labels (column A) replacing the old asset-class names
$ws.Range("A2").Value = "EUR001M Index"
$ws.Range("A3").Value = "LEF1TREU Index"
$ws.Range("A4").Value = "SX5R Index"
$ws.Range("A5").Value = "SXUSR Index`tUS"
$ws.Range("A6").Value = "BEGCGA Index"
$ws.Range("A7").Value = "LEC4TREU Index"
$ws.Range("A8").Value = "LEATTREU Index"

# Updated weights for "Opt Portfolio" (column B) and "Opt Portfolio with View" (column C)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.9999999999999998

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = [double]"3.05311331771918e-16"

$ws.Range("B4").Value = [double]"6.436803831546e-16"
$ws.Range("C4").Value = 0

$ws.Range("B5").Value = [double]"3.787206496253191e-16"
$ws.Range("C5").Value = [double]"3.200284730239067e-17"

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = [double]"5.560422236538279e-16"

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0

$ws.Range("B8").Value = [double]"2.517749162550946e-15"
$ws.Range("C8").Value = [double]"1.206421473722328e-15"
